$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of data rows (1-based table rows) to the five equations in that row.
# Row 1  -> table row 1
# Row 2  -> table row 5
# Row 3  -> table row 10
# Row 4  -> table row 15
# Row 5  -> table row 20
$updates = @{
    1  = @("93×91=8463", "46×39=1794", "31×50=1550", "88×69=6072", "41×26=1066")
    5  = @("27×18=486",  "25×54=1350", "25×43=1075", "18×84=1512", "96×30=2880")
    10 = @("89×94=8366", "80×32=2560", "17×19=323",  "55×57=3135", "58×52=3016")
    15 = @("14×79=1106", "46×42=1932", "25×73=1825", "11×23=253",  "26×82=2132")
    20 = @("62×24=1488", "12×66=792",  "96×80=7680", "76×98=7448", "86×99=8514")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $r = $cell.Range
        $r.End = $r.End - 1
        $r.Text = $values[$col - 1]
    }
}
